$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan")

# Update the "Project Lead:" cell to include the team members' names.
$ws.Range("A2").Value = "Project Lead: Brennan Sullivan and Jacob Coleman"

# Reset the view: scroll back to the top and select B2 (as if the user
# had just finished editing A2 and moved on).
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$ws.Range("B2").Select() | Out-Null
